# Simulate the two-stage buck (5V -> 20V -> 210V) "HV supply" concept.
# Adds a new plan column (H) to the boost-convertor calcs, alongside the
# existing B / D / E plans, and tweaks several of the existing plan inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HV supply")

# --- existing plan tweaks ---------------------------------------------
# Vin: plan B was 12V, plan D/E were 5V -> all three now feed from 20V
$ws.Range("B6").Value = 20
$ws.Range("D6").Value = 20
$ws.Range("E6").Value = 20

# switching frequency: 100kHz -> 47kHz for all three existing plans
$ws.Range("B7").Value = 47000
$ws.Range("D7").Value = 47000
$ws.Range("E7").Value = 47000

# plan D target Vout: 150 -> 120
$ws.Range("D14").Value = 120

# --- new plan: column H, "20V" (5V -> 20V pre-regulator stage) --------
$ws.Range("H3").Value = "20V"

$ws.Range("H5").Formula = "=6*(2.5+0.3)*0.001"

$ws.Range("H6").Value = 5

$ws.Range("H7").Value = 219000
$ws.Range("H7").NumberFormat = "0.00E+00"

$ws.Range("H8").Formula = "=1/H7"

$ws.Range("H9").Value = 0.000015
$ws.Range("H9").NumberFormat = "0.00E+00"

$ws.Range("H12").Formula = "=H5/H14/H6"

$ws.Range("H13").Formula = "=H8/H9*H6*H15"

$ws.Range("H14").Value = 20

$ws.Range("H15").Formula = "=1-(H6/H14)"
$ws.Range("H15").NumberFormat = "0.00%"

$ws.Range("H16").Formula = "=H8/(2*H9)*H14*H15*POWER(1-H15,2)"

$ws.Range("H17").Formula = "=H5/(1-H15)"

$ws.Range("H18").Formula = "=H5/(1-H15)+(H8/(2*H9)*H6*H15)"

$ws.Range("H19").Formula = "=H8/(2*H5)*POWER(H6/H14,2)*(H14-H6)"

# Conditional formatting on H9 (Iin-delta), mirroring the B9 / D9:E9 rules:
# green when <= IL-ripple threshold (B19), red when greater.
$fcLow = $ws.Range("H9").FormatConditions.Add(1, 8, '=$B$19')
$fcLow.Interior.Color = 13561798
$fcLow.Font.Color = 24832

$fcHigh = $ws.Range("H9").FormatConditions.Add(1, 5, '=$B$19')
$fcHigh.Interior.Color = 13551615
$fcHigh.Font.Color = 393372
